$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 215 (existing rows 215.. shift down to 217..),
# preserving formatting (date style on column D) from the row being pushed down.
$ws.Rows.Item(215).Resize(2).Insert()

# Row 215 - new record
$ws.Cells.Item(215, 1).Value = 7
$ws.Cells.Item(215, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(215, 3).Value = "Ñuble"
$ws.Cells.Item(215, 4).Value = 44636
$ws.Cells.Item(215, 5).Value = 16
$ws.Cells.Item(215, 6).Value = "Fruta"
$ws.Cells.Item(215, 7).Value = 100101
$ws.Cells.Item(215, 8).Value = "Berries"
$ws.Cells.Item(215, 9).Value = 100112025
$ws.Cells.Item(215, 10).Value = "Frutilla"
$ws.Cells.Item(215, 11).Value = "Sin especificar"
$ws.Cells.Item(215, 12).Value = "Primera"
$ws.Cells.Item(215, 13).Value = 120
$ws.Cells.Item(215, 14).Value = 6500
$ws.Cells.Item(215, 15).Value = 7000
$ws.Cells.Item(215, 16).Value = 6750
$ws.Cells.Item(215, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(215, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(215, 19).Value = 964
$ws.Cells.Item(215, 20).Value = 7

# Row 216 - new record
$ws.Cells.Item(216, 1).Value = 7
$ws.Cells.Item(216, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(216, 3).Value = "Ñuble"
$ws.Cells.Item(216, 4).Value = 44636
$ws.Cells.Item(216, 5).Value = 16
$ws.Cells.Item(216, 6).Value = "Fruta"
$ws.Cells.Item(216, 7).Value = 100101
$ws.Cells.Item(216, 8).Value = "Berries"
$ws.Cells.Item(216, 9).Value = 100112025
$ws.Cells.Item(216, 10).Value = "Frutilla"
$ws.Cells.Item(216, 11).Value = "Sin especificar"
$ws.Cells.Item(216, 12).Value = "Segunda"
$ws.Cells.Item(216, 13).Value = 60
$ws.Cells.Item(216, 14).Value = 6000
$ws.Cells.Item(216, 15).Value = 6000
$ws.Cells.Item(216, 16).Value = 6000
$ws.Cells.Item(216, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(216, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(216, 19).Value = 857
$ws.Cells.Item(216, 20).Value = 7
